$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows before row 393, pushing existing rows 393-431
# down to 395-433 (dimension grows from A1:T431 to A1:T433).
$ws.Range("A393:T393").EntireRow.Insert()
$ws.Range("A393:T393").EntireRow.Insert()

# New row 393: Naranja / Navel Late / Primera, fecha 2022-08-10 (44783)
$ws.Range("A393").Value = 7
$ws.Range("B393").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C393").Value = "Ñuble"
$ws.Range("D393").Value = 44783
$ws.Range("E393").Value = 16
$ws.Range("F393").Value = "Fruta"
$ws.Range("G393").Value = 100102
$ws.Range("H393").Value = "Cítricos"
$ws.Range("I393").Value = 100102005
$ws.Range("J393").Value = "Naranja"
$ws.Range("K393").Value = "Navel Late"
$ws.Range("L393").Value = "Primera"
$ws.Range("M393").Value = 120
$ws.Range("N393").Value = 5500
$ws.Range("O393").Value = 6000
$ws.Range("P393").Value = 5750
$ws.Range("Q393").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R393").Value = "Región de O'Higgins"
$ws.Range("S393").Value = 383
$ws.Range("T393").Value = 15

# New row 394: Naranja / Navel Late / Segunda, fecha 2022-08-10 (44783)
$ws.Range("A394").Value = 7
$ws.Range("B394").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C394").Value = "Ñuble"
$ws.Range("D394").Value = 44783
$ws.Range("E394").Value = 16
$ws.Range("F394").Value = "Fruta"
$ws.Range("G394").Value = 100102
$ws.Range("H394").Value = "Cítricos"
$ws.Range("I394").Value = 100102005
$ws.Range("J394").Value = "Naranja"
$ws.Range("K394").Value = "Navel Late"
$ws.Range("L394").Value = "Segunda"
$ws.Range("M394").Value = 60
$ws.Range("N394").Value = 5000
$ws.Range("O394").Value = 5000
$ws.Range("P394").Value = 5000
$ws.Range("Q394").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R394").Value = "Región de O'Higgins"
$ws.Range("S394").Value = 333
$ws.Range("T394").Value = 15

# Ensure the D393/D394 date cells carry the same number format as the other
# date cells in column D (custom "YYYY-MM-DD HH:MM:SS" format, style index 2).
$ws.Range("D393").NumberFormat = $ws.Range("D395").NumberFormat
$ws.Range("D394").NumberFormat = $ws.Range("D395").NumberFormat
